$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Line spacing: every paragraph goes from 288 (1.2x, "auto") to 360
#    (1.5x, "auto").  In the COM object model, LineSpacing is expressed
#    in points where 12pt == a single-spacing multiple of 1.0, so 1.5x
#    is 18 points while LineSpacingRule stays "multiple" (wdLineSpaceMultiple = 5).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 5
    $p.Format.LineSpacing = 18
}

# ---------------------------------------------------------------------
# 2. "My project has four source code modules..." paragraph gains a new
#    closing sentence explaining the reasoning behind the module split.
# ---------------------------------------------------------------------
$oldModules = "the module containing code for the linked list of books and the module containing code for an additional linked list I created to store rejected lines."
$newModules = $oldModules + "  I chose to create my modules this way because my main method was becoming lengthy, and I wanted to better organize my code."
$d.Content.Find.Execute($oldModules, $false, $false, $false, $false, $false, $true, 1, $false, $newModules, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Rewrite the "extra credit" paragraph's opening explanation.
# ---------------------------------------------------------------------
$oldExtra = "I think this project is worthy of extra credit because I attempted both of the commands on the PDF and they are both working.  I tested them by creating two new text files, books2.txt and books3.txt, and made minor changes to them to make sure they are working.  When I ran my program without the extra credit commands, I discovered they were coincidentally in alphabetical order by title and by author.  "
$newExtra = "I think this project is worthy of extra credit because I attempted both the –r command and the –a command on the PDF, and they are both working.  When I ran my program without the extra credit commands, I discovered they were coincidentally in alphabetical order by title and by author.  I tested the extra credit commands"
$d.Content.Find.Execute($oldExtra, $false, $false, $false, $false, $false, $true, 1, $false, $newExtra, 2) | Out-Null

# Locate the freshly written "...I tested the extra credit commands" anchor
# so we know exactly where the rest of the rewritten sentence (and the
# relocated _GoBack bookmark) belong.
$anchor = $d.Content
$anchor.Find.Execute("I tested the extra credit commands", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $anchor.End

$tailRange = $d.Range($splitPos, $splitPos)
$tailText = " by creating a series of new text files, copying and pasting the original data into them and making minor changes to the files, all to make sure the program is working.  "
$tailRange.InsertAfter($tailText)

# The _GoBack bookmark (a single, Word-maintained "last edit" marker) now
# belongs right after "...the extra credit commands" and before the
# newly-inserted tail text -- adding it here relocates it away from its
# old spot earlier in the document automatically.
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

Write-Host "done"
